# Trade #39 closed at 2026-02-16 21:29:06 - momentum DOWN +0.000%
# Append a new row (row 11) to the "momentum" worksheet describing this open trade.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("momentum")

$row = 11

$ws.Cells.Item($row, 1).Value = 39

# Column B holds a date-looking string ("2026-02-16"); force it to stay text
# (matching the other rows) instead of being auto-converted to a date serial.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "21:29:05"
$ws.Cells.Item($row, 4).Value = "momentum"
$ws.Cells.Item($row, 5).Value = "DOWN"
$ws.Cells.Item($row, 6).Value = 68644.355

# Column G (Exit Price) is blank for an open trade; use a quote-prefixed
# empty value so the cell is written out as an empty text cell rather than
# being omitted entirely.
$ws.Cells.Item($row, 7).Value = "'"
$ws.Cells.Item($row, 7).Style = "Normal"

$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.9
$ws.Cells.Item($row, 12).Value = "Downward momentum: -0.417% over 10 samples"

# Column M (Exit Reason) is blank for an open trade.
$ws.Cells.Item($row, 13).Value = "'"
$ws.Cells.Item($row, 13).Style = "Normal"

$ws.Cells.Item($row, 14).Value = 0
